$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one row per day, newest first, starting at row 2.
# This "auto update" run adds a new most-recent day on top, pushing all
# existing rows down by one (dimension grows by one row).

# Current newest date lives in A2 (as plain text, e.g. "2026-01-16").
$prevTopDateText = $ws.Range("A2").Value()
$prevTopDate = [DateTime]::ParseExact($prevTopDateText, "yyyy-MM-dd", $null)
$newDateText = $prevTopDate.AddDays(1).ToString("yyyy-MM-dd")

# Insert a fresh row above the current row 2; shifts rows 2..58 down to 3..59.
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting from the row above/below during
# Insert(); strip that so it matches the unstyled data rows below it.
$ws.Range("A2:D2").ClearFormats()

# Write the date as literal text (leading apostrophe forces text, avoiding
# Excel's automatic date-serial conversion), then reset the style so no
# stray number-format residue (quote-prefix style) is left behind.
$ws.Cells.Item(2, 1).Value = "'" + $newDateText
$ws.Range("A2").Style = "Normal"

# Prices are unchanged day over day in this feed.
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
